$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card1")
$ws.Rows.Item(3).Delete()
